{"js": "// Apply the three text edits described by the diff:\n// 1) Title: \" AUDIENCE METERING \" -> \" AUDIENCE MEASUREMENT \"\n// 2) \"Content consumption measurement...\" paragraph rewrite\n// 3) \"The classical view of audience metering...\" paragraph expansion\n\nconst body = context.document.body;\n\n// --- Change 1: title \"AUDIENCE METERING\" -> \"AUDIENCE MEASUREMENT\" ---\n{\n  const results = body.search(\" AUDIENCE METERING \", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Title text ' AUDIENCE METERING ' not found\");\n  }\n  results.items[0].insertText(\" AUDIENCE MEASUREMENT \", \"Replace\");\n  await context.sync();\n}\n\n// --- Change 2: \"Content consumption measurement ...\" intro paragraph ---\n{\n  const oldText =\n    \"Content consumption measurement is one of the aims the audience measurement technologies which include the analysis of users\\u2019 behaviour when consuming content or media services. Content consumption measurement is one of the main methods used by service providers or broadcasters, the main methods used by services providers or broadcasters to obtain useful data for refining service offerings or setting advertising rates. It\\u2019s applicability goes well beyond this. Without reliable audience data, many businesses will be reluctant to participate in the new delivery platforms.\";\n  const newText =\n    \"Content consumption measurement is one of the solutions the audience measurement technologies aim to bring about. Content consumption measurement is one of the main methods used by service providers or broadcasters to obtain useful data for refining service offerings or setting advertising rates. It\\u2019s applicability goes well beyond this. Without reliable audience data, many businesses will be reluctant to participate in the new delivery platforms.\";\n\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Paragraph 2 source text not found\");\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// --- Change 3: \"The classical view of audience metering ...\" definition sentence ---\n{\n  const oldText =\n    \" is the determination of the number of people who watched a particular TV program or channel, or listened to a certain radio station over a stated period of time. Direct and indirect methods of measurement are used\";\n  const newText =\n    \" is the estimation of the number of viewers who are tuned to a particular TV program or channel, or the number of listeners who are tuned to a particular radio program or channel. Audience measurement takes into account the behaviour of  the audience as well as their demographics. Direct and indirect methods of measurement are used\";\n\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Paragraph 3 source text not found\");\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the three text edits described by the diff:\n# 1) Title: \" AUDIENCE METERING \" -> \" AUDIENCE MEASUREMENT \"\n# 2) \"Content consumption measurement...\" paragraph rewrite\n# 3) \"The classical view of audience metering...\" paragraph expansion\n\n$d = $word.ActiveDocument\n\n# --- Change 1: title \"AUDIENCE METERING\" -> \"AUDIENCE MEASUREMENT\" ---\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Text = \" AUDIENCE METERING \"\n$rng1.Find.Forward = $true\n$rng1.Find.Wrap = 0\n$found1 = $rng1.Find.Execute()\nif ($found1) {\n    $rng1.Text = \" AUDIENCE MEASUREMENT \"\n}\n\n# --- Change 2: \"Content consumption measurement ...\" intro paragraph ---\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"Content consumption measurement is one of the aims the audience measurement technologies which include the analysis of users\u2019 behaviour when consuming content or media services. Content consumption measurement is one of the main methods used by service providers or broadcasters, the main methods used by services providers or broadcasters to obtain useful data for refining service offerings or setting advertising rates. It\u2019s applicability goes well beyond this. Without reliable audience data, many businesses will be reluctant to participate in the new delivery platforms.\"\n$rng2.Find.Forward = $true\n$rng2.Find.Wrap = 0\n$found2 = $rng2.Find.Execute()\nif ($found2) {\n    $rng2.Text = \"Content consumption measurement is one of the solutions the audience measurement technologies aim to bring about. Content consumption measurement is one of the main methods used by service providers or broadcasters to obtain useful data for refining service offerings or setting advertising rates. It\u2019s applicability goes well beyond this. Without reliable audience data, many businesses will be reluctant to participate in the new delivery platforms.\"\n}\n\n# --- Change 3: \"The classical view of audience metering ...\" definition sentence ---\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Text = \" is the determination of the number of people who watched a particular TV program or channel, or listened to a certain radio station over a stated period of time. Direct and indirect methods of measurement are used\"\n$rng3.Find.Forward = $true\n$rng3.Find.Wrap = 0\n$found3 = $rng3.Find.Execute()\nif ($found3) {\n    $rng3.Text = \" is the estimation of the number of viewers who are tuned to a particular TV program or channel, or the number of listeners who are tuned to a particular radio program or channel. Audience measurement takes into account the behaviour of  the audience as well as their demographics. Direct and indirect methods of measurement are used\"\n}\n\nWrite-Output \"found1:$found1 found2:$found2 found3:$found3\"\n"}
